$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '43.758.79'
Set-TextValue 'E2' '  +5.04%  '

Set-TextValue 'D3' '2.275.18'
Set-TextValue 'E3' '  +2.35%  '

Set-TextValue 'E4' '  +0.06%  '

Set-TextValue 'D5' '231.23'
Set-TextValue 'E5' '  +0.36%  '

Set-TextValue 'D6' '0.628'
Set-TextValue 'E6' '  +1.29%  '

Set-TextValue 'D7' '63.74'
Set-TextValue 'E7' '  +6.64%  '

Set-TextValue 'E8' '  +0.09%  '

Set-TextValue 'E9' '  +6.95%  '

Set-TextValue 'E10' '  +17.64%  '

Set-TextValue 'D11' '57.28'
Set-TextValue 'E11' '  -0.89%  '

Set-TextValue 'D12' '25.99'
Set-TextValue 'E12' '  +16.13%  '

Set-TextValue 'D13' '0.104'
Set-TextValue 'E13' '  +0.09%  '

Set-TextValue 'D14' '2.613.89'
Set-TextValue 'E14' '  +2.31%  '

Set-TextValue 'D15' '15.69'
Set-TextValue 'E15' '  +1.96%  '

Set-TextValue 'D16' '5.91'
Set-TextValue 'E16' '  +5.73%  '

Set-TextValue 'D17' '0.823'
Set-TextValue 'E17' '  +3.42%  '

Set-TextValue 'D18' '2.263.48'
Set-TextValue 'E18' '  +1.54%  '

Set-TextValue 'D19' '43.663.36'
Set-TextValue 'E19' '  +4.98%  '

Set-TextValue 'D20' '0.0000101'
Set-TextValue 'E20' '  +11.63%  '

Set-TextValue 'D21' '73.45'
Set-TextValue 'E21' '  +1.73%  '

Set-TextValue 'E22' '  +0.03%  '

Set-TextValue 'D23' '250.50'
Set-TextValue 'E23' '  +1.70%  '

Set-TextValue 'E24' '  +0.22%  '

Set-TextValue 'E25' '  +5.82%  '

Set-TextValue 'E26' '  -2.41%  '

Set-TextValue 'D27' '9.87'
Set-TextValue 'E27' '  +2.02%  '

Set-TextValue 'D28' '171.91'
Set-TextValue 'E28' '  +2.09%  '

Set-TextValue 'D29' '20.95'
Set-TextValue 'E29' '  +5.79%  '

Set-TextValue 'D30' '0.136'
Set-TextValue 'E30' '  -2.34%  '

Set-TextValue 'E31' '  +2.91%  '

Set-TextValue 'D32' '2.80'
Set-TextValue 'E32' '  +10.14%  '

Set-TextValue 'D33' '0.123'
Set-TextValue 'E33' '  +0.93%  '

Set-TextValue 'D34' '0.0686'
Set-TextValue 'E34' '  +5.57%  '

Set-TextValue 'D35' '5.05'
Set-TextValue 'E35' '  +1.83%  '

Set-TextValue 'D36' '4.74'
Set-TextValue 'E36' '  +1.80%  '

Set-TextValue 'B37' 'THORChain'
Set-TextValue 'C37' 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
Set-TextValue 'D37' '6.79'
Set-TextValue 'E37' '  +4.25%  '

Set-TextValue 'B38' 'RenderToken'
Set-TextValue 'C38' 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D38' '3.81'
Set-TextValue 'E38' '  +6.80%  '

Set-TextValue 'D39' '2.34'
Set-TextValue 'E39' '  -1.69%  '

Set-TextValue 'E40' '  +4.58%  '

Set-TextValue 'D41' '0.999'
Set-TextValue 'E41' '  -0.20%  '

Set-TextValue 'D42' '8.39'
Set-TextValue 'E42' '  -2.37%  '

Set-TextValue 'D43' '10.46'
Set-TextValue 'E43' '  +21.58%  '

Set-TextValue 'B44' 'InjectiveProtocol'
Set-TextValue 'C44' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue 'D44' '17.21'
Set-TextValue 'E44' '  +4.50%  '

Set-TextValue 'B45' 'Cronos'
Set-TextValue 'C45' 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue 'D45' '0.0962'
Set-TextValue 'E45' '  -0.33%  '

Set-TextValue 'D46' '1.21'
Set-TextValue 'E46' '  +0.20%  '

Set-TextValue 'D47' '97.59'
Set-TextValue 'E47' '  -0.70%  '

Set-TextValue 'D48' '4.42'
Set-TextValue 'E48' '  +0.47%  '

Set-TextValue 'D49' '1.478.35'
Set-TextValue 'E49' '  +0.42%  '

Set-TextValue 'D50' '2.34'
Set-TextValue 'E50' '  +4.87%  '

Set-TextValue 'D51' '0.000205'
Set-TextValue 'E51' '  -12.70%  '
